$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: clear the format-string that used to live in B3
$ws.Range("B3").Value = ""

# Row 4: used to hold "digits" numbers, now holds printf-style format strings
$ws.Range("B4").Value = "%m/%d/%Y %H:%M:%OS"
$ws.Range("C4").Value = "%5.1f"
$ws.Range("D4").Value = "%5.1f"
$ws.Range("E4").Value = "%.3f"
$ws.Range("F4").Value = "%.2f"
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = "%.2e"
$ws.Range("J4").Value = "%d"
$ws.Range("O4").Value = "%s"

# Update the saved selection to P7
$ws.Range("P7").Select()
